$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3551.4893
$ws.Range("I98").Value = 3340.6316
$ws.Range("J98").Value = 4441.778
$ws.Range("K98").Value = 3340.6316
$ws.Range("L98").Value = 4441.778
$ws.Range("M98").Value = -1842.6316
$ws.Range("N98").Value = -7437.778

$ws.Range("H111").Value = 9618785
$ws.Range("I111").Value = 15628045
$ws.Range("J111").Value = 3968.4
$ws.Range("K111").Value = 46884135
$ws.Range("L111").Value = 11905.2
$ws.Range("M111").Value = -46881068
$ws.Range("N111").Value = -18039.2

$ws.Range("H112").Value = 16899.908
$ws.Range("I112").Value = 1250
$ws.Range("J112").Value = 20377.666
$ws.Range("K112").Value = 3750
$ws.Range("L112").Value = 61132.99800000001
$ws.Range("M112").Value = -2642
$ws.Range("N112").Value = -63348.99800000001

$ws.Range("H122").Value = 3551.4893
$ws.Range("I122").Value = 3340.6316
$ws.Range("J122").Value = 4441.778
$ws.Range("K122").Value = 10021.8948
$ws.Range("L122").Value = 13325.334
$ws.Range("M122").Value = -7571.8948
$ws.Range("N122").Value = -18225.334

$ws.Range("H132").Value = 1429.7297
$ws.Range("I132").Value = 1403.0869
$ws.Range("K132").Value = 4209.2607
$ws.Range("M132").Value = -1679.2607

$ws.Range("H137").Value = 4269.069
$ws.Range("I137").Value = 5171.3076
$ws.Range("J137").Value = 3536
$ws.Range("K137").Value = 15513.9228
$ws.Range("L137").Value = 10608
$ws.Range("M137").Value = -12963.9228
$ws.Range("N137").Value = -15708

$ws.Range("H138").Value = 7154371
$ws.Range("I138").Value = 4544.3
$ws.Range("J138").Value = 25028938
$ws.Range("K138").Value = 13632.9
$ws.Range("L138").Value = 75086814
$ws.Range("M138").Value = -8492.900000000001
$ws.Range("N138").Value = -75097094

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1715932.2
$ws.Range("I32").Value = 1814507.9
$ws.Range("J32").Value = 15502.5
$ws.Range("K32").Value = 1814507.9
$ws.Range("L32").Value = 15502.5
$ws.Range("M32").Value = -1814220.9
$ws.Range("N32").Value = -16076.5

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 39318.668
$ws.Range("I74").Value = 48718.5
$ws.Range("K74").Value = 48718.5
$ws.Range("M74").Value = -47844.5

$ws.Range("H77").Value = 39318.668
$ws.Range("I77").Value = 48718.5
$ws.Range("K77").Value = 243592.5
$ws.Range("M77").Value = -239224.5

$ws.Range("H132").Value = 1230353.4
$ws.Range("I132").Value = 2110408.2
$ws.Range("J132").Value = 8055.0557
$ws.Range("K132").Value = 6331224.600000001
$ws.Range("L132").Value = 24165.1671
$ws.Range("M132").Value = -6328694.600000001
$ws.Range("N132").Value = -29225.1671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3340.7666
$ws.Range("I94").Value = 1747.7646
$ws.Range("J94").Value = 5423.923
$ws.Range("K94").Value = 1747.7646
$ws.Range("L94").Value = 5423.923
$ws.Range("M94").Value = -1296.7646
$ws.Range("N94").Value = -6325.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4865
$ws.Range("I16").Value = 2883.5557
$ws.Range("K16").Value = 2883.5557
$ws.Range("M16").Value = -2596.5557

$ws.Range("H31").Value = 13811.9375
$ws.Range("I31").Value = 4998.6
$ws.Range("J31").Value = 17818
$ws.Range("K31").Value = 4998.6
$ws.Range("L31").Value = 17818
$ws.Range("M31").Value = -4703.6
$ws.Range("N31").Value = -18408

$ws.Range("H34").Value = 13811.9375
$ws.Range("I34").Value = 4998.6
$ws.Range("J34").Value = 17818
$ws.Range("K34").Value = 4998.6
$ws.Range("L34").Value = 17818
$ws.Range("M34").Value = -4796.6
$ws.Range("N34").Value = -18222

$ws.Range("H86").Value = 18967098
$ws.Range("I86").Value = 7139097.5
$ws.Range("J86").Value = 47692240
$ws.Range("K86").Value = 7139097.5
$ws.Range("L86").Value = 47692240
$ws.Range("M86").Value = -7137974.5
$ws.Range("N86").Value = -47694486

$ws.Range("H89").Value = 18967098
$ws.Range("I89").Value = 7139097.5
$ws.Range("J89").Value = 47692240
$ws.Range("K89").Value = 35695487.5
$ws.Range("L89").Value = 238461200
$ws.Range("M89").Value = -35689871.5
$ws.Range("N89").Value = -238472432

$ws.Range("H99").Value = 6802.095
$ws.Range("I99").Value = 6152.5
$ws.Range("J99").Value = 7668.222
$ws.Range("K99").Value = 6152.5
$ws.Range("L99").Value = 7668.222
$ws.Range("M99").Value = -4654.5
$ws.Range("N99").Value = -10664.222

$ws.Range("H113").Value = 4865
$ws.Range("I113").Value = 2883.5557
$ws.Range("K113").Value = 2883.5557
$ws.Range("M113").Value = -713.5556999999999

$ws.Range("H126").Value = 6802.095
$ws.Range("I126").Value = 6152.5
$ws.Range("J126").Value = 7668.222
$ws.Range("K126").Value = 18457.5
$ws.Range("L126").Value = 23004.666
$ws.Range("M126").Value = -15987.5
$ws.Range("N126").Value = -27944.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4303.8047
$ws.Range("J68").Value = 4549.343
$ws.Range("L68").Value = 13648.029
$ws.Range("N68").Value = -15270.029

$ws.Range("H71").Value = 4303.8047
$ws.Range("J71").Value = 4549.343
$ws.Range("L71").Value = 40944.087
$ws.Range("N71").Value = -49056.087

$ws.Range("H97").Value = 199.5
$ws.Range("I97").Value = 199
$ws.Range("J97").Value = 199.66667
$ws.Range("K97").Value = 597
$ws.Range("L97").Value = 599.00001
$ws.Range("M97").Value = -101
$ws.Range("N97").Value = -1591.00001

$ws.Range("H119").Value = 4001
$ws.Range("I119").Value = 5005
$ws.Range("J119").Value = 2495
$ws.Range("K119").Value = 15015
$ws.Range("L119").Value = 7485
$ws.Range("M119").Value = -10177
$ws.Range("N119").Value = -17161

$ws.Range("H132").Value = 11925.643
$ws.Range("I132").Value = 6892.857
$ws.Range("J132").Value = 16958.428
$ws.Range("K132").Value = 62035.713
$ws.Range("L132").Value = 152625.852
$ws.Range("M132").Value = -59505.713
$ws.Range("N132").Value = -157685.852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6301.1113
$ws.Range("I132").Value = 3461.0833
$ws.Range("J132").Value = 11981.167
$ws.Range("K132").Value = 10383.2499
$ws.Range("L132").Value = 35943.501
$ws.Range("M132").Value = -7853.249899999999
$ws.Range("N132").Value = -41003.501

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4063.818
$ws.Range("I7").Value = 2428.2856
$ws.Range("J7").Value = 6926
$ws.Range("K7").Value = 2428.2856
$ws.Range("L7").Value = 6926
$ws.Range("M7").Value = -2316.2856
$ws.Range("N7").Value = -7150

$ws.Range("H126").Value = 4063.818
$ws.Range("I126").Value = 2428.2856
$ws.Range("J126").Value = 6926
$ws.Range("K126").Value = 7284.8568
$ws.Range("L126").Value = 20778
$ws.Range("M126").Value = -4814.8568
$ws.Range("N126").Value = -25718

$ws.Range("H136").Value = 9687.788
$ws.Range("I136").Value = 5460
$ws.Range("J136").Value = 11525.956
$ws.Range("K136").Value = 16380
$ws.Range("L136").Value = 34577.868
$ws.Range("M136").Value = -13830
$ws.Range("N136").Value = -39677.868

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 17775.467
$ws.Range("I113").Value = 28255.555
$ws.Range("J113").Value = 2055.3333
$ws.Range("K113").Value = 84766.66500000001
$ws.Range("L113").Value = 6165.999899999999
$ws.Range("M113").Value = -82596.66500000001
$ws.Range("N113").Value = -10505.9999

$ws.Range("H119").Value = 55653.332
$ws.Range("J119").Value = 55653.332
$ws.Range("L119").Value = 55653.332
$ws.Range("N119").Value = -65329.332

$ws.Range("H122").Value = 7638979
$ws.Range("I122").Value = 9884559
$ws.Range("K122").Value = 29653677
$ws.Range("M122").Value = -29651227

$ws.Range("H132").Value = 16706972
$ws.Range("I132").Value = 20018046
$ws.Range("J132").Value = 151601.8
$ws.Range("K132").Value = 60054138
$ws.Range("L132").Value = 454805.4
$ws.Range("M132").Value = -60051608
$ws.Range("N132").Value = -459865.4
